# Updated cryptos list on Mon Feb 27 20:57:53 UTC 2023 with GitHub Actions
# Refreshes the coin price/volume table (and the two-coin rank swaps) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "23.320.87"
$ws.Range("E2").Value = "  -1.01%  "
# Row 3
$ws.Range("D3").Value = "1.623.54"
$ws.Range("E3").Value = "  -0.88%  "
# Row 4
$ws.Range("E4").Value = "  +0.45%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.50%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.22%  "
# Row 7
$ws.Range("E7").Value = "  -0.08%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.21"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.71%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3611"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.85%  "
# Row 10
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.223"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.39%  "
# Row 11
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08073"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.61%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.54%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.40%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.544"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.02%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001244"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.55%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.210"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.72%  "
# Row 17
$ws.Range("D17").Value = "1.619.37"
$ws.Range("E17").Value = "  -0.94%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.60%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06911"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.57%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.40%  "
# Row 21
$ws.Range("E21").Value = "  +0.59%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.411"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.14%  "
# Row 23
$ws.Range("D23").Value = "23.316.48"
$ws.Range("E23").Value = "  -1.00%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.50%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.182"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.84%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.446"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.81%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.84%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.99%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.283"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.85%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.60%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.294"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.47%  "
# Row 32
$ws.Range("D32").Value = "1.802.14"
$ws.Range("E32").Value = "  -0.58%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.762"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.97%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.30%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9459"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.53%  "
# Row 36
$ws.Range("E36").Value = "  -2.30%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2515"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.29%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08821"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.55%  "
# Row 39
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.07107"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.27%  "
# Row 40
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.036"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.37%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.358"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.73%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7017"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.44%  "
# Row 43
$ws.Range("E43").Value = "  -1.82%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.23%  "
# Row 45
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6417"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.62%  "
# Row 46
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.54%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.308"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.77%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.987"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.25%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07972"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.12%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.196"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.07%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "125.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.48%  "
